$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3519
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 2797.5
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 2797.5
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -3449.5
$ws.Range("H33").Value = 26316620
$ws.Range("I33").Value = 605.08826
$ws.Range("K33").Value = 605.08826
$ws.Range("M33").Value = -376.08826
$ws.Range("H98").Value = 2533.75
$ws.Range("I98").Value = 1574
$ws.Range("K98").Value = 1574
$ws.Range("M98").Value = -76
$ws.Range("H111").Value = 1934.36
$ws.Range("I111").Value = 1329.1765
$ws.Range("J111").Value = 3220.375
$ws.Range("K111").Value = 3987.5295
$ws.Range("L111").Value = 9661.125
$ws.Range("M111").Value = -920.5295000000001
$ws.Range("N111").Value = -15795.125
$ws.Range("H113").Value = 4811
$ws.Range("I113").Value = 3700
$ws.Range("J113").Value = 5518
$ws.Range("K113").Value = 3700
$ws.Range("L113").Value = 5518
$ws.Range("M113").Value = -446
$ws.Range("N113").Value = -12026
$ws.Range("H116").Value = 55558276
$ws.Range("I116").Value = 2750
$ws.Range("J116").Value = 100002700
$ws.Range("K116").Value = 2750
$ws.Range("L116").Value = 100002700
$ws.Range("M116").Value = 692
$ws.Range("N116").Value = -100009584
$ws.Range("H122").Value = 2533.75
$ws.Range("I122").Value = 1574
$ws.Range("K122").Value = 4722
$ws.Range("M122").Value = -2272
$ws.Range("H132").Value = 1537.4667
$ws.Range("I132").Value = 1556
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 4668
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -2138
$ws.Range("N132").Value = -8060
$ws.Range("H138").Value = 5002840.5
$ws.Range("I138").Value = 1343.409
$ws.Range("J138").Value = 11115781
$ws.Range("K138").Value = 4030.227
$ws.Range("L138").Value = 33347343
$ws.Range("M138").Value = 1109.773
$ws.Range("N138").Value = -33357623

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 27326.666
$ws.Range("J24").Value = 27326.666
$ws.Range("L24").Value = 27326.666
$ws.Range("N24").Value = -28074.666
$ws.Range("H28").Value = 6571.5713
$ws.Range("I28").Value = 6571.5713
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 6571.5713
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -6379.5713
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 48000
$ws.Range("J31").Value = 48000
$ws.Range("L31").Value = 48000
$ws.Range("N31").Value = -48588
$ws.Range("H45").Value = 1955.3334
$ws.Range("I45").Value = 1639.6
$ws.Range("J45").Value = 2350
$ws.Range("K45").Value = 1639.6
$ws.Range("L45").Value = 2350
$ws.Range("M45").Value = -1262.6
$ws.Range("N45").Value = -3104
$ws.Range("H92").Value = 26404.666
$ws.Range("J92").Value = 26404.666
$ws.Range("L92").Value = 26404.666
$ws.Range("N92").Value = -31396.666
$ws.Range("H93").Value = 25005.6
$ws.Range("J93").Value = 25005.6
$ws.Range("L93").Value = 25005.6
$ws.Range("N93").Value = -29997.6
$ws.Range("H94").Value = 24990
$ws.Range("J94").Value = 24990
$ws.Range("L94").Value = 24990
$ws.Range("N94").Value = -26792
$ws.Range("H95").Value = 30208
$ws.Range("J95").Value = 30208
$ws.Range("L95").Value = 30208
$ws.Range("N95").Value = -35700
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H98").Value = 29271.834
$ws.Range("J98").Value = 29271.834
$ws.Range("L98").Value = 29271.834
$ws.Range("N98").Value = -35261.834
$ws.Range("H99").Value = 6571.5713
$ws.Range("I99").Value = 6571.5713
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6571.5713
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3576.5713
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 27326.666
$ws.Range("J100").Value = 27326.666
$ws.Range("L100").Value = 27326.666
$ws.Range("N100").Value = -29490.666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 11386.167
$ws.Range("J97").Value = 22796
$ws.Range("L97").Value = 22796
$ws.Range("N97").Value = -24778

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 987.7778
$ws.Range("I16").Value = 931.6667
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 931.6667
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -644.6667
$ws.Range("N16").Value = -1674
$ws.Range("H86").Value = 5948.8335
$ws.Range("I86").Value = 2870.2222
$ws.Range("J86").Value = 15184.667
$ws.Range("K86").Value = 2870.2222
$ws.Range("L86").Value = 15184.667
$ws.Range("M86").Value = -1747.2222
$ws.Range("N86").Value = -17430.667
$ws.Range("H89").Value = 5948.8335
$ws.Range("I89").Value = 2870.2222
$ws.Range("J89").Value = 15184.667
$ws.Range("K89").Value = 14351.111
$ws.Range("L89").Value = 75923.33499999999
$ws.Range("M89").Value = -8735.111000000001
$ws.Range("N89").Value = -87155.33499999999
$ws.Range("H94").Value = 2988.682
$ws.Range("I94").Value = 2850.818
$ws.Range("J94").Value = 3126.5454
$ws.Range("K94").Value = 2850.818
$ws.Range("L94").Value = 3126.5454
$ws.Range("M94").Value = -2399.818
$ws.Range("N94").Value = -4028.5454
$ws.Range("H99").Value = 1909
$ws.Range("I99").Value = 1868.6666
$ws.Range("J99").Value = 1939.25
$ws.Range("K99").Value = 1868.6666
$ws.Range("L99").Value = 1939.25
$ws.Range("M99").Value = -370.6666
$ws.Range("N99").Value = -4935.25
$ws.Range("H107").Value = 1686.1578
$ws.Range("I107").Value = 805.3333
$ws.Range("J107").Value = 2092.6924
$ws.Range("K107").Value = 805.3333
$ws.Range("L107").Value = 2092.6924
$ws.Range("M107").Value = 1114.6667
$ws.Range("N107").Value = -5932.6924
$ws.Range("H113").Value = 987.7778
$ws.Range("I113").Value = 931.6667
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 931.6667
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1238.3333
$ws.Range("N113").Value = -5440
$ws.Range("H122").Value = 950.875
$ws.Range("I122").Value = 799.4706
$ws.Range("J122").Value = 1318.5714
$ws.Range("K122").Value = 2398.4118
$ws.Range("L122").Value = 3955.7142
$ws.Range("M122").Value = 51.58820000000014
$ws.Range("N122").Value = -8855.7142
$ws.Range("H126").Value = 1909
$ws.Range("I126").Value = 1868.6666
$ws.Range("J126").Value = 1939.25
$ws.Range("K126").Value = 5605.9998
$ws.Range("L126").Value = 5817.75
$ws.Range("M126").Value = -3135.9998
$ws.Range("N126").Value = -10757.75
$ws.Range("H132").Value = 2535.4773
$ws.Range("I132").Value = 2084.257
$ws.Range("J132").Value = 4290.222
$ws.Range("K132").Value = 6252.771000000001
$ws.Range("L132").Value = 12870.666
$ws.Range("M132").Value = -3722.771000000001
$ws.Range("N132").Value = -17930.666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 83.083336
$ws.Range("I8").Value = 83.083336
$ws.Range("K8").Value = 249.250008
$ws.Range("M8").Value = -110.250008
$ws.Range("H113").Value = 571.8333
$ws.Range("I113").Value = 523.9091
$ws.Range("J113").Value = 647.1429000000001
$ws.Range("K113").Value = 1571.7273
$ws.Range("L113").Value = 1941.4287
$ws.Range("M113").Value = 598.2727
$ws.Range("N113").Value = -6281.4287
$ws.Range("H131").Value = 896.5599999999999
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 896.5599999999999
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2689.68
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12769.68

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2055.4211
$ws.Range("I102").Value = 2056.8667
$ws.Range("J102").Value = 2050
$ws.Range("K102").Value = 2056.8667
$ws.Range("L102").Value = 2050
$ws.Range("M102").Value = -434.8667
$ws.Range("N102").Value = -5294

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2845.3447
$ws.Range("I7").Value = 2160
$ws.Range("J7").Value = 2988.125
$ws.Range("K7").Value = 2160
$ws.Range("L7").Value = 2988.125
$ws.Range("M7").Value = -2048
$ws.Range("N7").Value = -3212.125
$ws.Range("H22").Value = 906.9231
$ws.Range("J22").Value = 899.1667
$ws.Range("L22").Value = 899.1667
$ws.Range("N22").Value = -1489.1667
$ws.Range("H27").Value = 906.9231
$ws.Range("J27").Value = 899.1667
$ws.Range("L27").Value = 899.1667
$ws.Range("N27").Value = -1113.1667
$ws.Range("H61").Value = 1972.6316
$ws.Range("I61").Value = 1567.6923
$ws.Range("K61").Value = 1567.6923
$ws.Range("M61").Value = -1365.6923
$ws.Range("H113").Value = 1972.6316
$ws.Range("I113").Value = 1567.6923
$ws.Range("K113").Value = 1567.6923
$ws.Range("M113").Value = 602.3077000000001
$ws.Range("H126").Value = 2845.3447
$ws.Range("I126").Value = 2160
$ws.Range("J126").Value = 2988.125
$ws.Range("K126").Value = 6480
$ws.Range("L126").Value = 8964.375
$ws.Range("M126").Value = -4010
$ws.Range("N126").Value = -13904.375
$ws.Range("H136").Value = 1736.909
$ws.Range("I136").Value = 1673.5667
$ws.Range("J136").Value = 2370.3333
$ws.Range("K136").Value = 5020.7001
$ws.Range("L136").Value = 7110.999899999999
$ws.Range("M136").Value = -2470.7001
$ws.Range("N136").Value = -12210.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 38211.484
$ws.Range("I122").Value = 1925.2222
$ws.Range("J122").Value = 81755
$ws.Range("K122").Value = 5775.6666
$ws.Range("L122").Value = 245265
$ws.Range("M122").Value = -3325.6666
$ws.Range("N122").Value = -250165
$ws.Range("H130").Value = 26166.666
$ws.Range("J130").Value = 26166.666
$ws.Range("L130").Value = 26166.666
$ws.Range("N130").Value = -36206.666
$ws.Range("H132").Value = 3319.743
$ws.Range("I132").Value = 3135.4
$ws.Range("J132").Value = 3780.6
$ws.Range("K132").Value = 9406.200000000001
$ws.Range("L132").Value = 11341.8
$ws.Range("M132").Value = -6876.200000000001
$ws.Range("N132").Value = -16401.8
